$wb = $excel.ActiveWorkbook
$src = $wb.Worksheets.Item("2021-Q4")
$newSheet = $wb.Worksheets.Add()
$src.Range("A1:H2").Copy($newSheet.Range("A1"))
$newSheet.Name = "2022-Q4"

$before = $wb.Worksheets.Item(2)
$newSheet.Move($before)

$ns2 = $wb.Worksheets.Item("2022-Q4")
$ns2.Range("D2").NumberFormat = "@"
$ns2.Range("D2").Value = "0.67"
$ns2.Range("E2").NumberFormat = "@"
$ns2.Range("E2").Value = "91.81"
$ns2.Range("F2").NumberFormat = "@"
$ns2.Range("F2").Value = "8.55"
$ns2.Range("G2").NumberFormat = "@"
$ns2.Range("G2").Value = "0.0573"

# Update summary (总计) sheet
$summary = $wb.Worksheets.Item(1)
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("B3").Value = "2021-Q4"
$summary.Range("D3").Value = 0.06

# New row 4, seeded from row 3's formatting (A3 carries the bold/border/center style)
$summary.Range("A3:D3").Copy($summary.Range("A4"))
$summary.Range("A4").Value = 2
$summary.Range("B4").Value = "2021-Q3"
$summary.Range("C4").Value = 1
$summary.Range("D4").Value = 0.05
